# Updated cryptos list on Fri Aug 16 04:40:59 UTC 2024 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto rows on the active sheet to the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row new values. $null means that column is unchanged for that row.
$updates = @(
    @{ Row = 2; D = "57.949.82"; E = "  -0.43%  " },
    @{ Row = 3; D = "2.571.29"; E = "  -2.66%  " },
    @{ Row = 4; D = $null; E = "  -0.01%  " },
    @{ Row = 5; D = "517.81"; E = "  -0.56%  " },
    @{ Row = 6; D = "142.25"; E = "  -1.16%  " },
    @{ Row = 7; D = $null; E = "  -0.27%  " },
    @{ Row = 8; D = "0.563"; E = "  -1.16%  " },
    @{ Row = 9; D = "2.587.67"; E = "  -2.16%  " },
    @{ Row = 10; D = "6.71"; E = "  +0.91%  " },
    @{ Row = 11; D = $null; E = "  -1.45%  " },
    @{ Row = 12; D = $null; E = "  -3.97%  " },
    @{ Row = 13; D = $null; E = "  -1.02%  " },
    @{ Row = 14; D = "3.026.84"; E = "  -2.62%  " },
    @{ Row = 15; D = "57.932.01"; E = "  -0.52%  " },
    @{ Row = 16; D = "20.29"; E = "  -2.64%  " },
    @{ Row = 17; D = $null; E = "  -1.77%  " },
    @{ Row = 18; D = "2.563.51"; E = "  -3.14%  " },
    @{ Row = 19; D = "340.83"; E = "  +0.95%  " },
    @{ Row = 21; D = $null; E = "  -1.99%  " },
    @{ Row = 22; D = "6.33"; E = "  +0.63%  " },
    @{ Row = 23; D = $null; E = "  -0.11%  " },
    @{ Row = 24; D = "65.31"; E = "  +1.16%  " },
    @{ Row = 25; D = $null; E = "  -1.57%  " },
    @{ Row = 26; D = $null; E = "  -5.34%  " },
    @{ Row = 27; D = $null; E = "  -0.25%  " },
    @{ Row = 28; D = "2.683.55"; E = "  -2.94%  " },
    @{ Row = 29; D = $null; E = "  -1.82%  " },
    @{ Row = 30; D = $null; E = "  -6.36%  " },
    @{ Row = 31; D = $null; E = "  -0.03%  " },
    @{ Row = 32; D = "6.21"; E = "  -6.77%  " },
    @{ Row = 33; D = "1.58"; E = "  -0.80%  " },
    @{ Row = 34; D = "18.67"; E = "  -0.87%  " },
    @{ Row = 35; D = "149.88"; E = "  -1.71%  " },
    @{ Row = 36; D = $null; E = "  -3.20%  " },
    @{ Row = 37; D = $null; E = "  -3.56%  " },
    @{ Row = 38; D = "0.866"; E = "  -5.19%  " },
    @{ Row = 39; D = "35.96"; E = "  -2.29%  " },
    @{ Row = 40; D = $null; E = "  -2.89%  " },
    @{ Row = 41; D = $null; E = "  +0.06%  " },
    @{ Row = 42; D = $null; E = "  -3.18%  " },
    @{ Row = 43; D = $null; E = "  -0.32%  " },
    @{ Row = 44; D = "269.66"; E = "  +0.13%  " },
    @{ Row = 45; D = "10.66"; E = "  +0.31%  " },
    @{ Row = 46; D = "0.0948"; E = "  -2.14%  " },
    @{ Row = 47; D = $null; E = "  -3.25%  " },
    @{ Row = 48; D = "18.80"; E = "  -3.05%  " },
    @{ Row = 49; D = $null; E = "  -3.28%  " },
    @{ Row = 50; D = $null; E = "  -0.04%  " },
    @{ Row = 51; D = "1.971.38"; E = "  -3.15%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($row, 4)
        # Some "Price" values look like plain numbers (e.g. "6.71", "0.0948").
        # Force the cell to Text format first so Excel keeps the exact
        # string instead of silently converting it to a numeric value.
        $trimmed = $u.D.Trim()
        $isNumericLooking = $trimmed -match '^[+-]?[0-9]*\.?[0-9]+$'
        if ($isNumericLooking) {
            $cellD.NumberFormat = "@"
        }
        $cellD.Value = $u.D
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
